$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: the "Openstaande bestelling/PR's" Whitelist marker (column D, "TRUE")
# had ended up on the wrong rows. Remove it from rows 13 and 15 ...
$ws.Range("D13").Value = $null
$ws.Range("D15").Value = $null

# ... and put it on the correct rows (36-39), copying an existing "TRUE" cell
# so the value reuses the shared string and keeps the same (default) style.
$ws.Range("D12").Copy($ws.Range("D36"))
$ws.Range("D12").Copy($ws.Range("D37"))
$ws.Range("D12").Copy($ws.Range("D38"))
$ws.Range("D12").Copy($ws.Range("D39"))

# Row 37 (Tjerk Jansen) was missing its e-mail address/hyperlink - add it,
# copying the same value/style used by the other rows' e-mail column.
$ws.Range("B36").Copy($ws.Range("B37"))
$ws.Hyperlinks.Add($ws.Range("B37"), "mailto:bram.gerrits@vhe.nl")
$ws.Range("B37").Style = "Hyperlink"

# Reflect where the user ended up working in the sheet.
$ws.Range("H36").Select()
